# Users sheet: add 4 new accounts (rows 23-26), per commit:
# "Se agrego ASUC26, Y CCD01 AL CC09" (new Sucursal/branch users)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Row 23: F00319 / Sucursal 319 (stored as a plain number, like the source data)
$ws.Range("A23").Value = "F00319"
$ws.Range("C23").Value = 319

# Row 24: F00044 / Sucursal "044" (kept as text, matching the existing
# zero-padded branch-code column formatting)
$ws.Range("A24").Value = "F00044"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").HorizontalAlignment = -4152
$ws.Range("C24").Value = "044"

# Row 25: OSANTOS / Casa central
$ws.Range("A25").Value = "OSANTOS"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").HorizontalAlignment = -4152
$ws.Range("C25").Value = "Casa central"

# Row 26: CRECERAD / Casa central
$ws.Range("A26").Value = "CRECERAD"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").HorizontalAlignment = -4152
$ws.Range("C26").Value = "Casa central"

# Update the sheet view: scroll so row 10 is at the top and select E28,
# matching the author's final cursor position when they finished editing.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E28").Select()
